$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.050.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.44%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.511.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.15%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'601.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.12%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'183.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +5.74%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.600"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.31%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +4.29%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -1.82%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.32%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'4.121.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.01%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'32.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +12.91%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +0.01%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'68.037.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.44%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +0.93%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.509.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E18").Value = "'  +1.35%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'14.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +4.30%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'398.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.81%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'8.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.45%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "'Litecoin"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "'73.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.44%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "'Polygon"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'0.546"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.20%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.38%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +2.80%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.16%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'10.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +5.09%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -0.70%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.992"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.51%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'6.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.44%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.40%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +0.52%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'24.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = "'  +1.24%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.08%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.68"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +2.28%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'164.52"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'1.97"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +3.01%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.880"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.56%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +4.13%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +6.57%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'4.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.79%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'27.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +3.02%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'26.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.27%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'Maker"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'2.855.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.09%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'Hedera"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'0.0741"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.27%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'42.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.26%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +0.77%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'348.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.51%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +0.02%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'34.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.91%  "
$ws.Range("E51").Style = "Normal"
